$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by exactly one day.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# Updated production values (column B) for rows 31-45, reflecting retrained model output.
$ws.Cells.Item(31, 2).Value = 6
$ws.Cells.Item(32, 2).Value = 17
$ws.Cells.Item(33, 2).Value = 32
$ws.Cells.Item(34, 2).Value = 50
$ws.Cells.Item(35, 2).Value = 67
$ws.Cells.Item(36, 2).Value = 85
$ws.Cells.Item(37, 2).Value = 95
$ws.Cells.Item(38, 2).Value = 109
$ws.Cells.Item(39, 2).Value = 129
$ws.Cells.Item(40, 2).Value = 146
$ws.Cells.Item(41, 2).Value = 150
$ws.Cells.Item(42, 2).Value = 170
$ws.Cells.Item(43, 2).Value = 184
$ws.Cells.Item(44, 2).Value = 192
$ws.Cells.Item(45, 2).Value = 199
